$d = $word.ActiveDocument
$t = $d.Tables(1)
$t.Cell(1,1).Range.Text = "533×9="
$t.Cell(1,2).Range.Text = "910×7="
$t.Cell(1,3).Range.Text = "472×7="
$t.Cell(1,4).Range.Text = "881×8="
$t.Cell(1,5).Range.Text = "836×7="
$t.Cell(5,1).Range.Text = "328×2="
$t.Cell(5,2).Range.Text = "650×2="
$t.Cell(5,3).Range.Text = "949×9="
$t.Cell(5,4).Range.Text = "215×2="
$t.Cell(5,5).Range.Text = "769×7="
$t.Cell(10,1).Range.Text = "534×3="
$t.Cell(10,2).Range.Text = "554×3="
$t.Cell(10,3).Range.Text = "590×6="
$t.Cell(10,4).Range.Text = "149×2="
$t.Cell(10,5).Range.Text = "506×5="
$t.Cell(15,1).Range.Text = "881×3="
$t.Cell(15,2).Range.Text = "778×6="
$t.Cell(15,3).Range.Text = "669×2="
$t.Cell(15,4).Range.Text = "969×8="
$t.Cell(15,5).Range.Text = "532×9="
$t.Cell(20,1).Range.Text = "817×3="
$t.Cell(20,2).Range.Text = "483×5="
$t.Cell(20,3).Range.Text = "137×8="
$t.Cell(20,4).Range.Text = "302×5="
$t.Cell(20,5).Range.Text = "666×8="
Write-Output "done"